# Workbook has two sheets: "Raw Data" (Raw Data responses) and "Data Insights"
# (computed summary / averages that reference "Raw Data" columns).
#
# The author's change inserts a brand-new "Additional Chars" data column into
# "Raw Data" right before the old column I ("Latency 1"), which pushes every
# column from I onward one slot to the right (I->J, J->K, ... N->O). Excel's
# own reference-tracking then keeps every formula on "Data Insights" that
# pointed at those shifted "Raw Data" columns pointing at the *same logical*
# column, i.e. each of those formulas' column reference moves one letter to
# the right as well.
#
# The commit also nudges the UI selection state around: the active sheet
# moves from "Data Insights" back to "Raw Data", with the cursor parked on
# the freshly inserted header cell, while "Data Insights" keeps its own
# selection synced to the (now shifted) first formula cell of the summary
# row.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Raw Data")
$ws2 = $wb.Worksheets.Item("Data Insights")

# Insert a new column at I (shifts old I:N -> J:O, and Excel auto-updates
# every formula reference elsewhere in the workbook that pointed into the
# shifted range, e.g. the AVERAGE(...) formulas on "Data Insights"!F21:J21).
$ws1.Columns("I:I").Insert()

# New column inherits the look of its neighbours (F:H, width 16.140625) --
# match that as closely as the exposed ColumnWidth property allows.
$ws1.Columns("I:I").ColumnWidth = 15.3

# Header text for the newly inserted column.
$ws1.Range("I4").Value = "Additional Chars"

# Restore per-sheet selections as left by the edit, and land back on
# "Raw Data" as the active tab.
$ws1.Range("I5").Select() | Out-Null
$ws2.Range("F21").Select() | Out-Null
$ws1.Activate() | Out-Null
